$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Method Description")
$ws.Rows.Item(3).Insert()
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)
$ws.Range("A3").Value = "TAB_TO_EXCEL_XML"
$ws.Range("B3").Value = "Static Method"
$ws.Range("C3").Value = "Public"
$ws.Range("D3").Value = "Any internal table to excel xml"
$excel.CutCopyMode = 0
[void]$ws.Range("D3").Select()
